$wb = $excel.ActiveWorkbook

# --- sheet1: move the stored selection from E10 to A16 (no longer the active tab) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("A16").Select()

# --- add the new "copy-to-verticalHeaderTableTest" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "copy-to-verticalHeaderTableTest"

# --- section title row ---
$a1 = $newSheet.Range("A1")
$a1.Value = "- vertical header table"
$titleSrc = $ws1.Range("A1")
$titleSrc.Copy()
$a1.PasteSpecial(-4122)   # xlPasteFormats

# --- vertical header column (C2:C4) ---
$c2 = $newSheet.Range("C2")
$c2.Value = "header1"
$c3 = $newSheet.Range("C3")
$c3.Value = "header2"
$c4 = $newSheet.Range("C4")
$c4.Value = "header3"

$headerSrc = $ws1.Range("A8")
$headerSrc.Copy()
$c2.PasteSpecial(-4122)  # xlPasteFormats
$c2.Interior.ThemeColor = 5   # theme="4" (Accent1) fill, matching the new blue header fill

$c2.Copy()
$c3.PasteSpecial(-4122)
$c4.PasteSpecial(-4122)

# --- make the new sheet the active tab ---
$newSheet.Select()
$newSheet.Range("A1").Select()

Write-Host "done"
